$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Remove the (empty placeholder) <w:pBdr/> from the "Figure 1" caption
#    paragraph (pStyle="normal3") by explicitly clearing all four border
#    sides - this drops the element entirely instead of writing w:val="none".
# ---------------------------------------------------------------------
$wdBorderTop = -1
$wdBorderLeft = -2
$wdBorderBottom = -3
$wdBorderRight = -4
$wdLineStyleNone = 0

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Style.NameLocal -eq "normal3") {
        $rng = $p.Range
        $rng.Borders.Item($wdBorderTop).LineStyle = $wdLineStyleNone
        $rng.Borders.Item($wdBorderLeft).LineStyle = $wdLineStyleNone
        $rng.Borders.Item($wdBorderBottom).LineStyle = $wdLineStyleNone
        $rng.Borders.Item($wdBorderRight).LineStyle = $wdLineStyleNone
        break
    }
}

# ---------------------------------------------------------------------
# 2. After the paragraph holding "{%chartImage5}", insert two new
#    paragraphs for "{%chartImage6}" and "{%chartImage7}", each split
#    into three runs: "{%chartImage", the digit, and "}".
#    Technique: build each chunk in its own freshly-inserted paragraph
#    (so Word doesn't merge the text into one run), then delete the
#    paragraph marks between them to splice the runs back into a single
#    paragraph while keeping them as separate <w:r> elements.
# ---------------------------------------------------------------------
function Insert-ChartImagePlaceholder {
    param($afterParaIndex, $digit)

    $anchor = $d.Paragraphs.Item($afterParaIndex)
    $anchor.Range.InsertParagraphAfter()

    $p1 = $d.Paragraphs.Item($afterParaIndex + 1)
    $r1 = $p1.Range
    $r1.Collapse(1)
    $r1.InsertAfter("{%chartImage")

    $p1b = $d.Paragraphs.Item($afterParaIndex + 1)
    $p1b.Range.InsertParagraphAfter()
    $p2 = $d.Paragraphs.Item($afterParaIndex + 2)
    $r2 = $p2.Range
    $r2.Collapse(1)
    $r2.InsertAfter($digit)

    $p2.Range.InsertParagraphAfter()
    $p3 = $d.Paragraphs.Item($afterParaIndex + 3)
    $r3 = $p3.Range
    $r3.Collapse(1)
    $r3.InsertAfter("}")

    # Merge paragraph 1 and 2 (delete the mark between them)
    $pa = $d.Paragraphs.Item($afterParaIndex + 1)
    $endA = $pa.Range.End
    $markA = $d.Range($endA - 1, $endA)
    $markA.Delete()

    # Merge the (now combined) paragraph with paragraph 3
    $pb = $d.Paragraphs.Item($afterParaIndex + 1)
    $endB = $pb.Range.End
    $markB = $d.Range($endB - 1, $endB)
    $markB.Delete()
}

$chartImage5Index = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.TrimEnd([char]13) -eq "{%chartImage5}") {
        $chartImage5Index = $i
        break
    }
}

Insert-ChartImagePlaceholder $chartImage5Index "6"
Insert-ChartImagePlaceholder ($chartImage5Index + 1) "7"

# ---------------------------------------------------------------------
# 3. Add <a:noFill/> to the chart picture's shape properties.
# ---------------------------------------------------------------------
for ($i = 1; $i -le $d.Shapes.Count; $i++) {
    $shp = $d.Shapes.Item($i)
    $shp.Fill.Visible = $false
}

# ---------------------------------------------------------------------
# 4. Lower-case the display names of the Heading1-6 and Caption styles
#    (style ids are unchanged).
# ---------------------------------------------------------------------
$renames = @{
    "Heading1" = "heading 1";
    "Heading2" = "heading 2";
    "Heading3" = "heading 3";
    "Heading4" = "heading 4";
    "Heading5" = "heading 5";
    "Heading6" = "heading 6";
    "Caption"  = "caption";
}
foreach ($styleId in $renames.Keys) {
    $st = $d.Styles.Item($styleId)
    $st.NameLocal = $renames[$styleId]
}

# ---------------------------------------------------------------------
# 5. Suppress automatic hyphenation on the "normal3" paragraph style.
# ---------------------------------------------------------------------
$normal3 = $d.Styles.Item("normal3")
$normal3.ParagraphFormat.Hyphenation = $false
